$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 91
$ws.Cells.Item(91,1).Value = 89
$ws.Cells.Item(91,2).Value = 6924568
$ws.Cells.Item(91,3).Value = "Mexico Liga de Expansion"
$ws.Cells.Item(91,4).Value = "Mexico Liga de Expansion"
$ws.Cells.Item(91,5).Value = 45214.92013888889
$ws.Cells.Item(91,6).Value = "Atletico Morelia"
$ws.Cells.Item(91,7).Value = "Atlante"
$ws.Cells.Item(91,8).Value = 0
$ws.Cells.Item(91,9).Value = 1
$ws.Cells.Item(91,10).Value = "A"
$ws.Cells.Item(91,11).Value = 2.4
$ws.Cells.Item(91,12).Value = 3
$ws.Cells.Item(91,13).Value = 2.875
$ws.Cells.Item(91,14).Value = 2.7
$ws.Cells.Item(91,15).Value = 3.1
$ws.Cells.Item(91,16).Value = 2.8
$ws.Cells.Item(91,17).Value = 0
$ws.Cells.Item(91,18).Value = 1.85
$ws.Cells.Item(91,19).Value = 1.95
$ws.Cells.Item(91,20).Value = 2.25
$ws.Cells.Item(91,21).Value = 1.975
$ws.Cells.Item(91,22).Value = 1.725
$ws.Cells.Item(91,23).Value = -1
$ws.Cells.Item(91,24).Value = -1
$ws.Cells.Item(91,25).Value = 1.8
$ws.Cells.Item(91,26).Value = -1
$ws.Cells.Item(91,27).Value = 0.95
$ws.Cells.Item(91,28).Value = -1
$ws.Cells.Item(91,29).Value = 0.7250000000000001

# Row 92
$ws.Cells.Item(92,1).Value = 90
$ws.Cells.Item(92,2).Value = 6924569
$ws.Cells.Item(92,3).Value = "Mexico Liga de Expansion"
$ws.Cells.Item(92,4).Value = "Mexico Liga de Expansion"
$ws.Cells.Item(92,5).Value = 45214.92013888889
$ws.Cells.Item(92,6).Value = "Venados FC"
$ws.Cells.Item(92,7).Value = "Dorados"
$ws.Cells.Item(92,8).Value = 4
$ws.Cells.Item(92,9).Value = 1
$ws.Cells.Item(92,10).Value = "H"
$ws.Cells.Item(92,11).Value = 1.615
$ws.Cells.Item(92,12).Value = 4
$ws.Cells.Item(92,13).Value = 4.5
$ws.Cells.Item(92,14).Value = 1.5
$ws.Cells.Item(92,15).Value = 4.75
$ws.Cells.Item(92,16).Value = 5.75
$ws.Cells.Item(92,17).Value = -1.25
$ws.Cells.Item(92,18).Value = 1.925
$ws.Cells.Item(92,19).Value = 1.875
$ws.Cells.Item(92,20).Value = 3
$ws.Cells.Item(92,21).Value = 1.75
$ws.Cells.Item(92,22).Value = 1.95
$ws.Cells.Item(92,23).Value = 0.5
$ws.Cells.Item(92,24).Value = -1
$ws.Cells.Item(92,25).Value = -1
$ws.Cells.Item(92,26).Value = 0.925
$ws.Cells.Item(92,27).Value = -1
$ws.Cells.Item(92,28).Value = 0.75
$ws.Cells.Item(92,29).Value = -1

# Row 186
$ws.Cells.Item(186,1).Value = 184
$ws.Cells.Item(186,2).Value = 7648957
$ws.Cells.Item(186,3).Value = "Mexico Liga de Expansion"
$ws.Cells.Item(186,4).Value = "Mexico Liga de Expansion"
$ws.Cells.Item(186,5).Value = 45361.5
$ws.Cells.Item(186,6).Value = "Unam Pumas U23"
$ws.Cells.Item(186,7).Value = "Tijuana U23"
$ws.Cells.Item(186,8).Value = 2
$ws.Cells.Item(186,9).Value = 0
$ws.Cells.Item(186,10).Value = "H"
$ws.Cells.Item(186,11).Value = 1.666
$ws.Cells.Item(186,12).Value = 3.5
$ws.Cells.Item(186,13).Value = 4.2
$ws.Cells.Item(186,14).Value = 1.533
$ws.Cells.Item(186,15).Value = 4.333
$ws.Cells.Item(186,16).Value = 6
$ws.Cells.Item(186,17).Value = -1.25
$ws.Cells.Item(186,18).Value = 2.025
$ws.Cells.Item(186,19).Value = 1.775
$ws.Cells.Item(186,20).Value = 2.75
$ws.Cells.Item(186,21).Value = 1.775
$ws.Cells.Item(186,22).Value = 2.025
$ws.Cells.Item(186,23).Value = 0.5329999999999999
$ws.Cells.Item(186,24).Value = -1
$ws.Cells.Item(186,25).Value = -1
$ws.Cells.Item(186,26).Value = 1.025
$ws.Cells.Item(186,27).Value = -1
$ws.Cells.Item(186,28).Value = -1
$ws.Cells.Item(186,29).Value = 1.025

# Row 187
$ws.Cells.Item(187,1).Value = 185
$ws.Cells.Item(187,2).Value = 7648958
$ws.Cells.Item(187,3).Value = "Mexico Liga de Expansion"
$ws.Cells.Item(187,4).Value = "Mexico Liga de Expansion"
$ws.Cells.Item(187,5).Value = 45361.5
$ws.Cells.Item(187,6).Value = "Monterrey U23"
$ws.Cells.Item(187,7).Value = "Mazatlan FC U23"
$ws.Cells.Item(187,8).Value = 4
$ws.Cells.Item(187,9).Value = 3
$ws.Cells.Item(187,10).Value = "H"
$ws.Cells.Item(187,11).Value = 2.375
$ws.Cells.Item(187,12).Value = 3.1
$ws.Cells.Item(187,13).Value = 2.75
$ws.Cells.Item(187,14).Value = 2.375
$ws.Cells.Item(187,15).Value = 3.4
$ws.Cells.Item(187,16).Value = 3
$ws.Cells.Item(187,17).Value = -0.25
$ws.Cells.Item(187,18).Value = 2
$ws.Cells.Item(187,19).Value = 1.8
$ws.Cells.Item(187,20).Value = 2.75
$ws.Cells.Item(187,21).Value = 1.95
$ws.Cells.Item(187,22).Value = 1.85
$ws.Cells.Item(187,23).Value = 1.375
$ws.Cells.Item(187,24).Value = -1
$ws.Cells.Item(187,25).Value = -1
$ws.Cells.Item(187,26).Value = 1
$ws.Cells.Item(187,27).Value = -1
$ws.Cells.Item(187,28).Value = 0.95
$ws.Cells.Item(187,29).Value = -1

# Row 227
$ws.Cells.Item(227,1).Value = 225
$ws.Cells.Item(227,2).Value = 7641725
$ws.Cells.Item(227,3).Value = "Mexico Liga de Expansion"
$ws.Cells.Item(227,4).Value = "Mexico Liga de Expansion"
$ws.Cells.Item(227,5).Value = 45392.92013888889
$ws.Cells.Item(227,6).Value = "Tlaxcala FC"
$ws.Cells.Item(227,7).Value = "Mineros de Zacatecas"
$ws.Cells.Item(227,8).Value = 0
$ws.Cells.Item(227,9).Value = 2
$ws.Cells.Item(227,10).Value = "A"
$ws.Cells.Item(227,11).Value = 2.75
$ws.Cells.Item(227,12).Value = 3.25
$ws.Cells.Item(227,13).Value = 2.25
$ws.Cells.Item(227,14).Value = 3.1
$ws.Cells.Item(227,15).Value = 3.6
$ws.Cells.Item(227,16).Value = 2.15
$ws.Cells.Item(227,17).Value = 0.25
$ws.Cells.Item(227,18).Value = 1.9
$ws.Cells.Item(227,19).Value = 1.9
$ws.Cells.Item(227,20).Value = 3
$ws.Cells.Item(227,21).Value = 1.925
$ws.Cells.Item(227,22).Value = 1.875
$ws.Cells.Item(227,23).Value = -1
$ws.Cells.Item(227,24).Value = -1
$ws.Cells.Item(227,25).Value = 1.15
$ws.Cells.Item(227,26).Value = -1
$ws.Cells.Item(227,27).Value = 0.8999999999999999
$ws.Cells.Item(227,28).Value = -1
$ws.Cells.Item(227,29).Value = 0.875

# Row 228
$ws.Cells.Item(228,1).Value = 226
$ws.Cells.Item(228,2).Value = 7641726
$ws.Cells.Item(228,3).Value = "Mexico Liga de Expansion"
$ws.Cells.Item(228,4).Value = "Mexico Liga de Expansion"
$ws.Cells.Item(228,5).Value = 45393.00347222222
$ws.Cells.Item(228,6).Value = "Universidad Guadalajara"
$ws.Cells.Item(228,7).Value = "Tepatitlan FC"
$ws.Cells.Item(228,8).Value = 3
$ws.Cells.Item(228,9).Value = 2
$ws.Cells.Item(228,10).Value = "H"
$ws.Cells.Item(228,11).Value = 1.25
$ws.Cells.Item(228,12).Value = 5.5
$ws.Cells.Item(228,13).Value = 7.5
$ws.Cells.Item(228,14).Value = 1.3
$ws.Cells.Item(228,15).Value = 5.75
$ws.Cells.Item(228,16).Value = 7.5
$ws.Cells.Item(228,17).Value = -1.5
$ws.Cells.Item(228,18).Value = 1.775
$ws.Cells.Item(228,19).Value = 1.925
$ws.Cells.Item(228,20).Value = 3
$ws.Cells.Item(228,21).Value = 1.975
$ws.Cells.Item(228,22).Value = 1.825
$ws.Cells.Item(228,23).Value = 0.3
$ws.Cells.Item(228,24).Value = -1
$ws.Cells.Item(228,25).Value = -1
$ws.Cells.Item(228,26).Value = -1
$ws.Cells.Item(228,27).Value = 0.925
$ws.Cells.Item(228,28).Value = 0.9750000000000001
$ws.Cells.Item(228,29).Value = -1

# Row 229 (new row, copy number formats for A and E from row 228)
$ws.Cells.Item(228,1).Copy() | Out-Null
$ws.Cells.Item(229,1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(228,5).Copy() | Out-Null
$ws.Cells.Item(229,5).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Cells.Item(229,1).Value = 227
$ws.Cells.Item(229,2).Value = 7641728
$ws.Cells.Item(229,3).Value = "Mexico Liga de Expansion"
$ws.Cells.Item(229,4).Value = "Mexico Liga de Expansion"
$ws.Cells.Item(229,5).Value = 45395.83680555555
$ws.Cells.Item(229,6).Value = "Oaxaca"
$ws.Cells.Item(229,7).Value = "Atletico Morelia"
$ws.Cells.Item(229,11).Value = 2.25
$ws.Cells.Item(229,12).Value = 3.5
$ws.Cells.Item(229,13).Value = 2.7
$ws.Cells.Item(229,14).Value = 2.25
$ws.Cells.Item(229,15).Value = 3.5
$ws.Cells.Item(229,16).Value = 2.7
$ws.Cells.Item(229,17).Value = -0.25
$ws.Cells.Item(229,18).Value = 2
$ws.Cells.Item(229,19).Value = 1.8
$ws.Cells.Item(229,20).Value = 2.75
$ws.Cells.Item(229,21).Value = 1.975
$ws.Cells.Item(229,22).Value = 1.825
$ws.Cells.Item(229,23).Value = 0
$ws.Cells.Item(229,24).Value = 0
$ws.Cells.Item(229,25).Value = 0
$ws.Cells.Item(229,26).Value = 0
$ws.Cells.Item(229,27).Value = 0
